# Weekly data refresh: insert two new price-record rows for the most
# recent reporting week right after the current last "Primera/Segunda"
# pair block (row 454), shifting all subsequent rows down by two.
#
# This mirrors the source data feed's pattern: each insert adds a
# "$/caja 36 atados" row followed by its paired "$/docena de atados"
# row for the same market date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 455:456 - everything from the old row 455
# down shifts to 457 onward (old 530:531 become 532:533).
$ws.Rows("455:456").Insert()

# New row 455 - "$/caja 36 atados" record
$ws.Range("A455").Value = 6
$ws.Range("B455").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C455").Value = "Metropolitana"
$ws.Range("D455").Value = 44476
$ws.Range("E455").Value = 13
$ws.Range("F455").Value = 100112040
$ws.Range("G455").Value = "Cilantro"
$ws.Range("H455").Value = "Sin especificar"
$ws.Range("I455").Value = "Primera"
$ws.Range("J455").Value = 720
$ws.Range("K455").Value = 4000
$ws.Range("L455").Value = 4500
$ws.Range("M455").Value = 4243
$ws.Range("N455").Value = "`$/caja 36 atados"
$ws.Range("O455").Value = "Región Metropolitana"
$ws.Range("P455").Value = 118
$ws.Range("Q455").Value = 36
$ws.Range("R455").Value = "Hortaliza"

# New row 456 - "$/docena de atados" record (same market date)
$ws.Range("A456").Value = 6
$ws.Range("B456").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C456").Value = "Metropolitana"
$ws.Range("D456").Value = 44476
$ws.Range("E456").Value = 13
$ws.Range("F456").Value = 100112040
$ws.Range("G456").Value = "Cilantro"
$ws.Range("H456").Value = "Sin especificar"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 470
$ws.Range("K456").Value = 6500
$ws.Range("L456").Value = 7000
$ws.Range("M456").Value = 6702
$ws.Range("N456").Value = "`$/docena de atados"
$ws.Range("O456").Value = "Región Metropolitana"
$ws.Range("P456").Value = 2234
$ws.Range("Q456").Value = 3
$ws.Range("R456").Value = "Hortaliza"
